# Aggiornato funzionamento della priorita' (ordine invertito):
# le tre commesse nelle righe 2, 4 e 5 (stessa macchina R10) vengono
# riordinate - la commessa che prima era pianificata per ultima (riga 5)
# passa in testa (riga 2), e le altre due scalano di conseguenza - e gli
# orari di inizio/fine setup e lavorazione vengono ricalcolati in base
# al nuovo ordine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Riga 2: ora ospita la commessa che prima era in riga 5 (251308)
$ws.Range("A2").Value = 251308
$ws.Range("D2").Value = 81.9672131147541
$ws.Range("H2").Value = "2025-04-24 08:41:58"
$ws.Range("I2").Value = 5000

# Riga 4: ora ospita la commessa che prima era in riga 2 (251168)
$ws.Range("A4").Value = 251168
$ws.Range("D4").Value = 142.3114754098361
$ws.Range("E4").Value = "2025-04-24 08:41:58"
$ws.Range("F4").Value = "2025-04-24 09:01:58"
$ws.Range("G4").Value = "2025-04-24 09:01:58"
$ws.Range("H4").Value = "2025-04-24 11:24:16"
$ws.Range("I4").Value = 8681

# Riga 5: ora ospita la commessa che prima era in riga 4 (251167)
$ws.Range("A5").Value = 251167
$ws.Range("D5").Value = 173.655737704918
$ws.Range("E5").Value = "2025-04-24 11:24:16"
$ws.Range("F5").Value = "2025-04-24 11:44:16"
$ws.Range("G5").Value = "2025-04-24 11:44:16"
$ws.Range("I5").Value = 10593
